# Weekly fruit/vegetable price update: a new record is inserted at the top
# of the data block (row 342), pushing the existing historical rows down by
# one (342-408 -> 343-409).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 342, shifting rows 342:408 down to 343:409
$ws.Rows.Item(342).Insert()

# Populate the newly inserted row 342 with the new weekly record
$ws.Range("A342").Value = 9
$ws.Range("B342").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C342").Value = "Metropolitana"
$ws.Range("D342").Value = 45015
$ws.Range("E342").Value = 13
$ws.Range("F342").Value = 300000001
$ws.Range("G342").Value = "Rabanito"
$ws.Range("H342").Value = "Sin especificar"
$ws.Range("I342").Value = "Primera"
$ws.Range("J342").Value = 7000
$ws.Range("K342").Value = 3000
$ws.Range("L342").Value = 3000
$ws.Range("M342").Value = 3000
$ws.Range("N342").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O342").Value = "Provincia de Chacabuco"
$ws.Range("P342").Value = 30
$ws.Range("Q342").Value = 100
$ws.Range("R342").Value = "Hortaliza"
